# Add new data row for "Число дет. муз, худ школ - musartschool (шт.) (8017010)"
# to the "All features" sheet, mirroring the existing library/cultureorg rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New lookup-table entry (columns H:J, row 35) -------------------------
# Set values in the same order the source file lists the new shared
# strings (long description, short name, date range) so the rebuilt
# sharedStrings table lines up with the target.
$ws.Range("C36").Value = "Число дет. муз, худ школ - musartschool (шт.) (8017010)"
$ws.Range("H35").Value = "musartschool"
$ws.Range("J35").Value = "2007 - 2017"
$ws.Range("I35").Value = 58942

# Match formatting of the existing rows directly above the new cells.
$ws.Range("H34").Copy()
$ws.Range("H35").PasteSpecial(-4122)

$ws.Range("I34").Copy()
$ws.Range("I35").PasteSpecial(-4122)

$ws.Range("J34").Copy()
$ws.Range("J35").PasteSpecial(-4122)

# --- New blank row 37 (keeps the trailing empty styled row) ---------------
$ws.Range("C35").Copy()
$ws.Range("C37").PasteSpecial(-4122)

# --- Cosmetic sheet-level tweaks from the diff -----------------------------
$ws.Columns.Item(3).ColumnWidth = 52.5
$ws.Range("D40").Select()

$excel.CutCopyMode = 0
